$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "53-43=10"
$t.Cell(1, 2).Range.Text = "76+5=81"
$t.Cell(1, 3).Range.Text = "19+49=68"
$t.Cell(1, 4).Range.Text = "58-53=5"
$t.Cell(1, 5).Range.Text = "3+37=40"
$t.Cell(2, 1).Range.Text = "11+54=65"
$t.Cell(2, 2).Range.Text = "97-96=1"
$t.Cell(2, 3).Range.Text = "33+19=52"
$t.Cell(2, 4).Range.Text = "2+86=88"
$t.Cell(2, 5).Range.Text = "87-16=71"
$t.Cell(3, 1).Range.Text = "22+12=34"
$t.Cell(3, 2).Range.Text = "49+15=64"
$t.Cell(3, 3).Range.Text = "95-77=18"
$t.Cell(3, 4).Range.Text = "16+4=20"
$t.Cell(3, 5).Range.Text = "46-25=21"
$t.Cell(4, 1).Range.Text = "31+1=32"
$t.Cell(4, 2).Range.Text = "60-4=56"
$t.Cell(4, 3).Range.Text = "55-17=38"
$t.Cell(4, 4).Range.Text = "37+2=39"
$t.Cell(4, 5).Range.Text = "25+41=66"
$t.Cell(5, 1).Range.Text = "60+18=78"
$t.Cell(5, 2).Range.Text = "67-4=63"
$t.Cell(5, 3).Range.Text = "79-27=52"
$t.Cell(5, 4).Range.Text = "80-42=38"
$t.Cell(5, 5).Range.Text = "48-14=34"
$t.Cell(6, 1).Range.Text = "78-8=70"
$t.Cell(6, 2).Range.Text = "93-88=5"
$t.Cell(6, 3).Range.Text = "64-42=22"
$t.Cell(6, 4).Range.Text = "54+44=98"
$t.Cell(6, 5).Range.Text = "30+37=67"
$t.Cell(7, 1).Range.Text = "89-18=71"
$t.Cell(7, 2).Range.Text = "87-22=65"
$t.Cell(7, 3).Range.Text = "99-39=60"
$t.Cell(7, 4).Range.Text = "94-81=13"
$t.Cell(7, 5).Range.Text = "12+72=84"
$t.Cell(8, 1).Range.Text = "2+35=37"
$t.Cell(8, 2).Range.Text = "67-3=64"
$t.Cell(8, 3).Range.Text = "69+24=93"
$t.Cell(8, 4).Range.Text = "76-19=57"
$t.Cell(8, 5).Range.Text = "8+14=22"
$t.Cell(9, 1).Range.Text = "11+47=58"
$t.Cell(9, 2).Range.Text = "62-21=41"
$t.Cell(9, 3).Range.Text = "50+24=74"
$t.Cell(9, 4).Range.Text = "45+37=82"
$t.Cell(9, 5).Range.Text = "72+2=74"
$t.Cell(10, 1).Range.Text = "58+15=73"
$t.Cell(10, 2).Range.Text = "13+83=96"
$t.Cell(10, 3).Range.Text = "24+56=80"
$t.Cell(10, 4).Range.Text = "25+32=57"
$t.Cell(10, 5).Range.Text = "91+8=99"
$t.Cell(11, 1).Range.Text = "29+28=57"
$t.Cell(11, 2).Range.Text = "64+6=70"
$t.Cell(11, 3).Range.Text = "88-44=44"
$t.Cell(11, 4).Range.Text = "23+74=97"
$t.Cell(11, 5).Range.Text = "81-5=76"
$t.Cell(12, 1).Range.Text = "53+38=91"
$t.Cell(12, 2).Range.Text = "17+80=97"
$t.Cell(12, 3).Range.Text = "19+27=46"
$t.Cell(12, 4).Range.Text = "39+17=56"
$t.Cell(12, 5).Range.Text = "7+42=49"
$t.Cell(13, 1).Range.Text = "43+36=79"
$t.Cell(13, 2).Range.Text = "47-9=38"
$t.Cell(13, 3).Range.Text = "66-48=18"
$t.Cell(13, 4).Range.Text = "80-60=20"
$t.Cell(13, 5).Range.Text = "6+72=78"
$t.Cell(14, 1).Range.Text = "41+31=72"
$t.Cell(14, 2).Range.Text = "61-49=12"
$t.Cell(14, 3).Range.Text = "43-10=33"
$t.Cell(14, 4).Range.Text = "15+55=70"
$t.Cell(14, 5).Range.Text = "36-34=2"
$t.Cell(15, 1).Range.Text = "63-43=20"
$t.Cell(15, 2).Range.Text = "53-49=4"
$t.Cell(15, 3).Range.Text = "34-27=7"
$t.Cell(15, 4).Range.Text = "54+16=70"
$t.Cell(15, 5).Range.Text = "95-7=88"
$t.Cell(16, 1).Range.Text = "4+17=21"
$t.Cell(16, 2).Range.Text = "2+47=49"
$t.Cell(16, 3).Range.Text = "46-15=31"
$t.Cell(16, 4).Range.Text = "45+33=78"
$t.Cell(16, 5).Range.Text = "60+9=69"
$t.Cell(17, 1).Range.Text = "22+53=75"
$t.Cell(17, 2).Range.Text = "9+72=81"
$t.Cell(17, 3).Range.Text = "47+38=85"
$t.Cell(17, 4).Range.Text = "63-48=15"
$t.Cell(17, 5).Range.Text = "43+44=87"
$t.Cell(18, 1).Range.Text = "84-16=68"
$t.Cell(18, 2).Range.Text = "68+30=98"
$t.Cell(18, 3).Range.Text = "38+12=50"
$t.Cell(18, 4).Range.Text = "34+65=99"
$t.Cell(18, 5).Range.Text = "71-66=5"
$t.Cell(19, 1).Range.Text = "68+24=92"
$t.Cell(19, 2).Range.Text = "14+54=68"
$t.Cell(19, 3).Range.Text = "85-69=16"
$t.Cell(19, 4).Range.Text = "90-6=84"
$t.Cell(19, 5).Range.Text = "59+14=73"
$t.Cell(20, 1).Range.Text = "15+6=21"
$t.Cell(20, 2).Range.Text = "16+81=97"
$t.Cell(20, 3).Range.Text = "59+4=63"
$t.Cell(20, 4).Range.Text = "7+1=8"
$t.Cell(20, 5).Range.Text = "23+18=41"
